# Scheduled-runner update: refresh computed Profit columns (H, I, J, K, L, M, N)
# across the ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR sheets.
$wb = $excel.ActiveWorkbook

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

# Row 12
$ws.Range("H12").Value = 50099.5
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

# Row 51
$ws.Range("H51").Value = 14393.458
$ws.Range("J51").Value = 14384.454
$ws.Range("L51").Value = 14384.454
$ws.Range("N51").Value = -15352.454

# Row 107
$ws.Range("H107").Value = 918.2222
$ws.Range("I107").Value = 854.58826
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 854.58826
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = 1065.41174
$ws.Range("N107").Value = -5840

# Row 137
$ws.Range("H137").Value = 1972.2858
$ws.Range("I137").Value = 1272.375
$ws.Range("J137").Value = 2905.5
$ws.Range("K137").Value = 3817.125
$ws.Range("L137").Value = 8716.5
$ws.Range("M137").Value = -1267.125
$ws.Range("N137").Value = -13816.5

# Row 138
$ws.Range("H138").Value = 5672.4473
$ws.Range("I138").Value = 3812.5
$ws.Range("J138").Value = 5891.2646
$ws.Range("K138").Value = 11437.5
$ws.Range("L138").Value = 17673.7938
$ws.Range("M138").Value = -6297.5
$ws.Range("N138").Value = -27953.7938

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

# Row 45
$ws.Range("H45").Value = 2298.5
$ws.Range("I45").Value = 2298
$ws.Range("K45").Value = 2298
$ws.Range("M45").Value = -1921

# Row 61
$ws.Range("H61").Value = 43481130
$ws.Range("I61").Value = 47621736
$ws.Range("K61").Value = 47621736
$ws.Range("M61").Value = -47621524

# Row 136
$ws.Range("H136").Value = 43481130
$ws.Range("I136").Value = 47621736
$ws.Range("K136").Value = 142865208
$ws.Range("M136").Value = -142862658

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

# Row 107
$ws.Range("H107").Value = 91703.37
$ws.Range("I107").Value = 592.25
$ws.Range("K107").Value = 592.25
$ws.Range("M107").Value = 1327.75

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 15983.19
$ws.Range("I31").Value = 13427.375
$ws.Range("J31").Value = 17556
$ws.Range("K31").Value = 13427.375
$ws.Range("L31").Value = 17556
$ws.Range("M31").Value = -13132.375
$ws.Range("N31").Value = -18146

# Row 34
$ws.Range("H34").Value = 15983.19
$ws.Range("I34").Value = 13427.375
$ws.Range("J34").Value = 17556
$ws.Range("K34").Value = 13427.375
$ws.Range("L34").Value = 17556
$ws.Range("M34").Value = -13225.375
$ws.Range("N34").Value = -17960

# Row 132
$ws.Range("H132").Value = 142860590
$ws.Range("J132").Value = 4000
$ws.Range("L132").Value = 12000
$ws.Range("N132").Value = -17060

# Row 134
$ws.Range("H134").Value = 31252804
$ws.Range("I134").Value = 41669240
$ws.Range("K134").Value = 125007720
$ws.Range("M134").Value = -125005185

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

# Row 34
$ws.Range("H34").Value = 970.8
$ws.Range("J34").Value = 1601.4
$ws.Range("L34").Value = 4804.200000000001
$ws.Range("N34").Value = -4972.200000000001

# Row 39
$ws.Range("H39").Value = 2999.3333
$ws.Range("I39").Value = 999
$ws.Range("J39").Value = 7000
$ws.Range("K39").Value = 2997
$ws.Range("L39").Value = 21000
$ws.Range("M39").Value = -2703
$ws.Range("N39").Value = -21588

# Row 40
$ws.Range("H40").Value = 91.09999999999999
$ws.Range("J40").Value = 500
$ws.Range("L40").Value = 2000
$ws.Range("N40").Value = -2138

# Row 55
$ws.Range("H55").Value = 948.75
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 1131.6666
$ws.Range("K55").Value = 1200
$ws.Range("L55").Value = 3394.9998
$ws.Range("M55").Value = -1023
$ws.Range("N55").Value = -3748.9998

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")

# Row 97
$ws.Range("H97").Value = 1051.5
$ws.Range("I97").Value = 589.75
$ws.Range("K97").Value = 589.75
$ws.Range("M97").Value = -93.75

# Row 113
$ws.Range("H113").Value = 73714.86
$ws.Range("I113").Value = 92459.91
$ws.Range("K113").Value = 92459.91
$ws.Range("M113").Value = -90289.91

# Row 132
$ws.Range("H132").Value = 4631511.5
$ws.Range("I132").Value = 4631511.5
$ws.Range("K132").Value = 13894534.5
$ws.Range("M132").Value = -13892004.5

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 2036.3334
$ws.Range("I22").Value = 2118.2856
$ws.Range("J22").Value = 1749.5
$ws.Range("K22").Value = 2118.2856
$ws.Range("L22").Value = 1749.5
$ws.Range("M22").Value = -1823.2856
$ws.Range("N22").Value = -2339.5

# Row 27
$ws.Range("H27").Value = 2036.3334
$ws.Range("I27").Value = 2118.2856
$ws.Range("J27").Value = 1749.5
$ws.Range("K27").Value = 2118.2856
$ws.Range("L27").Value = 1749.5
$ws.Range("M27").Value = -2011.2856
$ws.Range("N27").Value = -1963.5

# Row 50
$ws.Range("H50").Value = 50084
$ws.Range("J50").Value = 50084
$ws.Range("L50").Value = 50084
$ws.Range("N50").Value = -51358

# Row 61
$ws.Range("H61").Value = 1420
$ws.Range("I61").Value = 1420
$ws.Range("K61").Value = 1420
$ws.Range("M61").Value = -1218

# Row 113
$ws.Range("H113").Value = 1420
$ws.Range("I113").Value = 1420
$ws.Range("K113").Value = 1420
$ws.Range("M113").Value = 750

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")

# Row 62
$ws.Range("H62").Value = 4795.3335
$ws.Range("I62").Value = 4698
$ws.Range("K62").Value = 4698
$ws.Range("M62").Value = -4074

# Row 65
$ws.Range("H65").Value = 4795.3335
$ws.Range("I65").Value = 4698
$ws.Range("K65").Value = 23490
$ws.Range("M65").Value = -20370

# Row 81
$ws.Range("H81").Value = 2055.3
$ws.Range("I81").Value = 2161.4443
$ws.Range("J81").Value = 1100
$ws.Range("K81").Value = 4322.8886
$ws.Range("L81").Value = 2200
$ws.Range("M81").Value = -3261.8886
$ws.Range("N81").Value = -4322

# Row 84
$ws.Range("H84").Value = 2055.3
$ws.Range("I84").Value = 2161.4443
$ws.Range("J84").Value = 1100
$ws.Range("K84").Value = 21614.443
$ws.Range("L84").Value = 11000
$ws.Range("M84").Value = -16310.443
$ws.Range("N84").Value = -21608

# Row 100 (L100 takes the value that used to live in M100; M100 is removed)
$ws.Range("H100").Value = 1432.6666
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 1432.6666
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2865.3332
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -3947.3332

# Row 132
$ws.Range("H132").Value = 14290727
$ws.Range("I132").Value = 19232894
$ws.Range("K132").Value = 57698682
$ws.Range("M132").Value = -57696152

# Row 136
$ws.Range("H136").Value = 12197476
$ws.Range("I136").Value = 12502329
$ws.Range("K136").Value = 37506987
$ws.Range("M136").Value = -37504437
